# ---------------------------------------------------------------------------
# Implements: "Implement AddRemoveUserOrganizationUnits and related workflows"
#
# 1. Populates the "取得" (Get) sheet's Table14 with sample user rows that
#    exercise organization-unit add/remove scenarios.
# 2. Adds a brand-new "組織単位の追加・削除" (Add/Remove organization units)
#    worksheet at the end of the workbook, with its own table (Table136) and
#    one sample/result row, mirroring the other request/result sheets
#    (作成, 編集, 削除).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "取得" sheet — fill in the sample data rows (2-11) of Table14.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("取得")

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "admin"
$ws.Range("E2").Value = "admin@defaulttenant.com"
$ws.Range("F2").Value = "Default, Test"
$ws.Range("G2").Value = "Administrator, AssetsManager"
$ws.Range("H2").Value = "Active"
$ws.Rows.Item(2).RowHeight = 29

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "HostRobot"
$ws.Range("F3").Value = "Default"
$ws.Range("G3").Value = "Robot"
$ws.Range("H3").Value = "Active"

$ws.Range("A4").Value = 15
$ws.Range("B4").Value = "I"
$ws.Range("F4").Value = "Default"
$ws.Range("G4").Value = "Robot"
$ws.Range("H4").Value = "Active"

$ws.Range("A5").Value = 16
$ws.Range("B5").Value = "J"
$ws.Range("F5").Value = "Default"
$ws.Range("G5").Value = "Robot"
$ws.Range("H5").Value = "Active"

$ws.Range("A6").Value = 20
$ws.Range("B6").Value = "x"
$ws.Range("C6").Value = "x"
$ws.Range("H6").Value = "Active"

$ws.Range("A7").Value = 31
$ws.Range("B7").Value = "A"
$ws.Range("C7").Value = "A"
$ws.Range("D7").Value = "A"
$ws.Range("E7").Value = "a@a.com"
$ws.Range("H7").Value = "Active"

$ws.Range("A8").Value = 32
$ws.Range("B8").Value = "B"
$ws.Range("C8").Value = "B"
$ws.Range("D8").Value = "B"
$ws.Range("E8").Value = "b@b.com"
$ws.Range("F8").Value = "Test"
$ws.Range("H8").Value = "Active"

$ws.Range("A9").Value = 33
$ws.Range("B9").Value = "C"
$ws.Range("C9").Value = "C"
$ws.Range("D9").Value = "C"
$ws.Range("E9").Value = "c@c.com"
$ws.Range("G9").Value = "Administrator"
$ws.Range("H9").Value = "Active"

$ws.Range("A10").Value = 35
$ws.Range("B10").Value = "D"
$ws.Range("C10").Value = "D"
$ws.Range("D10").Value = "D"
$ws.Range("E10").Value = "d@d.com"
$ws.Range("F10").Value = "Default, Test"
$ws.Range("G10").Value = "Administrator, MachineViewer"
$ws.Range("H10").Value = "Active"
$ws.Rows.Item(10).RowHeight = 29

$ws.Range("A11").Value = 36
$ws.Range("B11").Value = "E"
$ws.Range("C11").Value = "E"
$ws.Range("D11").Value = "E"
$ws.Range("E11").Value = "e@e.com"
$ws.Range("F11").Value = "Test"
$ws.Range("G11").Value = "MachineViewer"
$ws.Range("H11").Value = "Active"

# Selection moved to B11 as part of the edit.
$ws.Range("B11").Select()

# ---------------------------------------------------------------------------
# 2) New "組織単位の追加・削除" sheet at the end of the workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "組織単位の追加・削除"

$newSheet.Range("A1").Value = "ユーザー名"
$newSheet.Range("B1").Value = "追加される組織単位名"
$newSheet.Range("C1").Value = "削除される組織単位名"
$newSheet.Range("D1").Value = "結果"

$newSheet.Range("A2").Value = "E"
$newSheet.Range("B2").Value = "Test"
$newSheet.Range("C2").Value = "Default"
$newSheet.Range("D2").Value = "成功"

$newSheet.Columns.Item(1).ColumnWidth = 18.45
$newSheet.Columns.Item(2).ColumnWidth = 24.27
$newSheet.Columns.Item(3).ColumnWidth = 24.27

$tbl = $newSheet.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $newSheet.Range("A1:D101"), [System.Reflection.Missing]::Value, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table136"
$tbl.TableStyle = "TableStyleMedium3"

$newSheet.Range("A2").Select()

Write-Output "done"
